$wb = $excel.ActiveWorkbook

# --- Sheet "CONT.SES" (sheet1) ---
$wsCont = $wb.Worksheets.Item("CONT.SES")

# Swap the color values for two of the registries (Azul <-> Preta)
$wsCont.Range("B8").Value = "Preta"
$wsCont.Range("B10").Value = "Azul"

# New column F width
$wsCont.Columns.Item(6).ColumnWidth = 17.83

# Add criteria cells for the COUNTIFS example with two criteria
$wsCont.Range("F2").Value = "Preta"
$wsCont.Range("F2").HorizontalAlignment = -4108
$wsCont.Range("F2").VerticalAlignment = -4108

$wsCont.Range("G2").Value = "Feminino"

$wsCont.Range("F3").Formula = "=COUNTIFS(B2:B14,F2,C2:C14,G2)"
$wsCont.Range("F3").HorizontalAlignment = -4108

# Data validation lists
$wsCont.Range("F2").Validation.Add(3, 1, 1, """Amarela,Azul,Preta""")
$wsCont.Range("H2").Validation.Add(3, 1, 1, """Masculino,Feminino""")
$wsCont.Range("G2").Validation.Add(3, 1, 1, """Masculino,Feminino""")

$wsCont.Activate()
$wsCont.Range("F6").Select()

# --- Sheet "Pratica" (sheet2) ---
$wsPratica = $wb.Worksheets.Item("Pratica")

$wsPratica.Range("H4").Formula = "=COUNTIFS(B:B,""Coordenador"")"
$wsPratica.Range("H8").Formula = "=COUNTIFS(E:E,"">31/12/2018"")"

$wsPratica.Activate()
$wsPratica.Range("H8").Select()
